# Updated cryptos list on Wed Sep 25 06:57:17 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.939.70'

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.623.17'
$ws.Range("E3").Value = '  -0.64%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  +0.00%  '

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.29'

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.19'
$ws.Range("E6").Value = '  +2.67%  '

# Row 7 - USDC
$ws.Range("E7").Value = '  -0.05%  '

# Row 8 - XRP
$ws.Range("E8").Value = '  +0.36%  '

# Row 9 - Dogecoin
$ws.Range("E9").Value = '  +1.21%  '

# Row 10 - Cardano
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.382'
$ws.Range("E10").Value = '  +4.96%  '

# Row 11 - Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.60'
$ws.Range("E11").Value = '  -0.03%  '

# Row 12 - TRON
$ws.Range("E12").Value = '  -0.89%  '

# Row 13 - Avalanche
$ws.Range("E13").Value = '  +1.13%  '

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.092.95'
$ws.Range("E14").Value = '  -0.74%  '

# Row 15 - WrappedBTC
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '63.794.53'
$ws.Range("E15").Value = '  +0.94%  '

# Row 16 - ShibaInu
$ws.Range("E16").Value = '  +2.82%  '

# Row 17 - WrappedEther
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.618.19'
$ws.Range("E17").Value = '  -0.96%  '

# Row 18 - Chainlink
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.17'
$ws.Range("E18").Value = '  +6.68%  '

# Row 19 - Polkadot
$ws.Range("E19").Value = '  +2.93%  '

# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '349.56'
$ws.Range("E20").Value = '  +2.43%  '

# Row 21 - Uniswap
$ws.Range("E21").Value = '  -0.05%  '

# Row 22 - Dai
$ws.Range("E22").Value = '  -0.05%  '

# Row 23 - LEO
$ws.Range("E23").Value = '  +2.09%  '

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.26'
$ws.Range("E24").Value = '  -0.24%  '

# Row 25 - SuiNetwork
$ws.Range("E25").Value = '  +13.87%  '

# Row 26 - Fetch.AI
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.68'
$ws.Range("E26").Value = '  +2.12%  '

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.20'
$ws.Range("E27").Value = '  +5.06%  '

# Row 28 - Aptos
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.11'
$ws.Range("E28").Value = '  +3.37%  '

# Row 29 - Kaspa
$ws.Range("E29").Value = '  +0.37%  '

# Row 30 - Bittensor
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '547.48'
$ws.Range("E30").Value = '  -0.40%  '

# Row 31 - Binance-PegBSC-USD
$ws.Range("E31").Value = '  -0.03%  '

# Row 32 - PancakeSwap
$ws.Range("E32").Value = '  +0.33%  '

# Row 33 - PEPE
$ws.Range("E33").Value = '  +5.52%  '

# Row 34 - ImmutableX
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.75'
$ws.Range("E34").Value = '  +0.28%  '

# Row 35 - NEARProtocol
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.23'
$ws.Range("E35").Value = '  -0.19%  '

# Row 36 - Monero
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '167.76'
$ws.Range("E36").Value = '  +0.86%  '

# Row 37 - PolygonEcosystemToken
$ws.Range("E37").Value = '  +0.39%  '

# Row 38 - FirstDigitalUSD
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.02%  '

# Row 39 - Stacks
$ws.Range("E39").Value = '  +4.07%  '

# Row 40 - EthereumClassic
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.40'
$ws.Range("E40").Value = '  +2.05%  '

# Row 41 - USDe
$ws.Range("E41").Value = '  +0.09%  '

# Row 42 - Aave
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '168.19'
$ws.Range("E42").Value = '  -0.05%  '

# Row 43 - OKB
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.88'
$ws.Range("E43").Value = '  +0.65%  '

# Row 44 - Filecoin
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.91'
$ws.Range("E44").Value = '  +4.66%  '

# Row 45 - Hedera
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0589'
$ws.Range("E45").Value = '  +2.36%  '

# Row 46 - InjectiveProtocol
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.45'
$ws.Range("E46").Value = '  -4.53%  '

# Row 47 - Mantle
$ws.Range("E47").Value = '  +0.63%  '

# Row 48 - was VeChain, now dogwifhat (rows 48/49 order changed with new values)
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.99'
$ws.Range("E48").Value = '  +12.36%  '

# Row 49 - was dogwifhat, now VeChain
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0245'
$ws.Range("E49").Value = '  +0.54%  '

# Row 50 - Stellar
$ws.Range("E50").Value = '  +0.53%  '

# Row 51 - EnergySwap
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.19'
$ws.Range("E51").Value = '  +2.36%  '

Write-Host "Applied cryptos update"
